$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.129.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.106.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "349.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5171"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4456"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.72"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08960"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.174"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.72"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.110.75"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.227"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.733"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001148"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.007"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.82"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +7.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06692"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.244"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.225.76"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.86"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.349"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.359.99"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.95"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.543"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.31"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.71"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.174"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.628"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.257"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.979"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.45"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.923"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02579"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06830"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2310"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.64"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6826"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.287"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.26"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.319"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6391"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000365"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.655"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.222"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.86"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07252"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.82%  "
